$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 309, shifting rows 309:350 down to 310:351
$ws.Rows.Item(309).Insert(-4121)  # -4121 = xlShiftDown

# Populate the newly inserted row 309 with the new record's data
$ws.Cells.Item(309, 1).Value = 8
$ws.Cells.Item(309, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(309, 3).Value = "Coquimbo"
$ws.Cells.Item(309, 4).Value = [DateTime]"2022-12-23"
$ws.Cells.Item(309, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(309, 5).Value = 4
$ws.Cells.Item(309, 6).Value = 100112012
$ws.Cells.Item(309, 7).Value = "Espinaca"
$ws.Cells.Item(309, 8).Value = "Sin especificar"
$ws.Cells.Item(309, 9).Value = "Primera"
$ws.Cells.Item(309, 10).Value = 2000
$ws.Cells.Item(309, 11).Value = 500
$ws.Cells.Item(309, 12).Value = 600
$ws.Cells.Item(309, 13).Value = 550
$ws.Cells.Item(309, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(309, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(309, 16).Value = 1100
$ws.Cells.Item(309, 17).Value = 0.5
$ws.Cells.Item(309, 18).Value = "Hortaliza"
